# Update "想去人数" (want-to-go count) values in column F on the
# "展览" (Exhibition) and "全部类型" (All types) worksheets, as produced by a
# newer scrape of the source site (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 769
$ws1.Range("F4").Value  = 1543
$ws1.Range("F6").Value  = 101
$ws1.Range("F7").Value  = 168
$ws1.Range("F8").Value  = 6329
$ws1.Range("F12").Value = 5413
$ws1.Range("F14").Value = 181
$ws1.Range("F15").Value = 1205
$ws1.Range("F21").Value = 310
$ws1.Range("F24").Value = 3876
$ws1.Range("F25").Value = 166

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 769
$ws4.Range("F5").Value  = 1543
$ws4.Range("F7").Value  = 101
$ws4.Range("F8").Value  = 168
$ws4.Range("F9").Value  = 6329
$ws4.Range("F13").Value = 5413
$ws4.Range("F15").Value = 181
$ws4.Range("F16").Value = 1205
$ws4.Range("F22").Value = 310
$ws4.Range("F25").Value = 3876
$ws4.Range("F27").Value = 166
